$d = $word.ActiveDocument

# --- First paragraph: add a paragraph border (5pt space on each edge), ---
# --- widen the left indent, drop the trailing space-only run, and      ---
# --- update the placeholder text                                      ---
$p1 = $d.Paragraphs(1)

# Add paragraph border (top/left/bottom/right), each with w:space="5"
$p1.Range.Borders.DistanceFromTop = 5
$p1.Range.Borders.DistanceFromLeft = 5
$p1.Range.Borders.DistanceFromBottom = 5
$p1.Range.Borders.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# The paragraph currently ends with "...ID** " (a trailing run holding a
# single space) followed by the paragraph mark. Drop that trailing space
# character -- this removes the now-empty trailing run entirely.
$paraEnd = $p1.Range.End
$trailingSpace = $d.Range($paraEnd - 2, $paraEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# Update the placeholder text carried by the (now sole) run in the paragraph
[void]$d.Content.Find.Execute("**ID__AFFARS_pgi_5319_topic_2__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5319__ID**", 2)
